$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.007.00'
$ws.Range('E2').Value = '  -0.25%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.551.85'
$ws.Range('E3').Value = '  -0.09%  '

# Row 4
$ws.Range('E4').Value = '  -0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.78'
$ws.Range('E5').Value = '  -2.09%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '196.86'
$ws.Range('E6').Value = '  +5.90%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'

# Row 8
$ws.Range('E8').Value = '  -0.07%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.211'
$ws.Range('E9').Value = '  -1.86%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.658'
$ws.Range('E10').Value = '  +0.58%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.14'
$ws.Range('E11').Value = '  +0.27%  '

# Row 12
$ws.Range('E12').Value = '  -1.96%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.58'
$ws.Range('E13').Value = '  +0.33%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.100.81'
$ws.Range('E14').Value = '  -0.49%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '601.76'
$ws.Range('E15').Value = '  -5.32%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.20'
$ws.Range('E16').Value = '  +1.21%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.201.51'
$ws.Range('E17').Value = '  -0.05%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.72'
$ws.Range('E18').Value = '  -2.00%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.537.61'
$ws.Range('E19').Value = '  -0.87%  '

# Row 20
$ws.Range('E20').Value = '  +0.62%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.17'
$ws.Range('E22').Value = '  +3.39%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.32'
$ws.Range('E23').Value = '  +7.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.60'
$ws.Range('E24').Value = '  +0.38%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.63'
$ws.Range('E25').Value = '  -2.30%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.13'
$ws.Range('E26').Value = '  +2.80%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('E27').Value = '  -0.99%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.71'
$ws.Range('E28').Value = '  +1.80%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.69'
$ws.Range('E29').Value = '  -3.09%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.50'
$ws.Range('E30').Value = '  +21.70%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.13'
$ws.Range('E31').Value = '  +0.79%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.77'
$ws.Range('E32').Value = '  +3.63%  '

# Row 33
$ws.Range('E33').Value = '  +0.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.49'
$ws.Range('E34').Value = '  -0.63%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0834'
$ws.Range('E35').Value = '  +6.71%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.746.98'
$ws.Range('E36').Value = '  +6.38%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.11'
$ws.Range('E37').Value = '  -4.85%  '

# Row 39
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.63'
$ws.Range('E39').Value = '  +2.74%  '

# Row 40
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.396'
$ws.Range('E40').Value = '  -1.75%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.95'
$ws.Range('E41').Value = '  -1.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '500.75'
$ws.Range('E42').Value = '  -5.90%  '

# Row 43
$ws.Range('E43').Value = '  -1.22%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0458'
$ws.Range('E44').Value = '  -0.64%  '

# Row 45
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  -3.85%  '

# Row 46
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.140'
$ws.Range('E46').Value = '  -2.52%  '

# Row 47
$ws.Range('E47').Value = '  -0.99%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.24%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.70'
$ws.Range('E49').Value = '  -5.59%  '

# Row 50
$ws.Range('E50').Value = '  +1.70%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.30'
$ws.Range('E51').Value = '  +11.45%  '
